# GBDS OCTOBER 2025 FILES | UPDATED FILES
# Fill in purchase detail rows 19-21 on the "PE, SEPTEMBER" sheet with
# October 2025 invoice data (date, invoice/OR number, and gross amount),
# which cascades through the dependent formulas in columns K, L and M
# and the totals row 32.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PE, SEPTEMBER")

# Row 19
$ws.Range("C19").Value = 45961
$ws.Range("G19").Value = 518062443
$ws.Range("I19").Formula = "=326410-12465.8"

# Row 20
$ws.Range("C20").Value = 45961
$ws.Range("G20").Value = 518062417
$ws.Range("I20").Formula = "=294720-31077.6"

# Row 21
$ws.Range("C21").Value = 45961
$ws.Range("G21").Value = 518061959
$ws.Range("I21").Formula = "=1193940-49636.8"

# Update the active selection to match the author's final cursor position
$ws.Range("I22").Select()
